# Auto-generated edit script: apply Kujata_Profits market-price refresh values
# (values sourced from the commit diff; columns H-N per row on sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Cells.Item(11, 8).Value = 141  # H11: 171.25 -> 141
$ws.Cells.Item(11, 9).Value = 141  # I11: 171.25 -> 141
$ws.Cells.Item(11, 11).Value = 141  # K11: 171.25 -> 141
$ws.Cells.Item(11, 13).Value = -1  # M11: -31.25 -> -1
# Row 21
$ws.Cells.Item(21, 8).Value = 23611.4  # H21: 29919 -> 23611.4
$ws.Cells.Item(21, 9).Value = 18500  # I21: 0 -> 18500
$ws.Cells.Item(21, 10).Value = 27019  # J21: 29919 -> 27019
$ws.Cells.Item(21, 11).Value = 18500  # K21: 0 -> 18500
$ws.Cells.Item(21, 12).Value = 27019  # L21: 29919 -> 27019
$ws.Cells.Item(21, 13).Value = -18032  # M21: <MISSING> -> -18032
$ws.Cells.Item(21, 14).Value = -27955  # N21: -30855 -> -27955
# Row 23
$ws.Cells.Item(23, 8).Value = 23611.4  # H23: 29919 -> 23611.4
$ws.Cells.Item(23, 9).Value = 18500  # I23: 0 -> 18500
$ws.Cells.Item(23, 10).Value = 27019  # J23: 29919 -> 27019
$ws.Cells.Item(23, 11).Value = 18500  # K23: 0 -> 18500
$ws.Cells.Item(23, 12).Value = 27019  # L23: 29919 -> 27019
$ws.Cells.Item(23, 13).Value = -18266  # M23: <MISSING> -> -18266
$ws.Cells.Item(23, 14).Value = -27487  # N23: -30387 -> -27487
# Row 32
$ws.Cells.Item(32, 8).Value = 1833  # H32: 1212.625 -> 1833
$ws.Cells.Item(32, 9).Value = 499  # I32: 760.4 -> 499
$ws.Cells.Item(32, 10).Value = 2500  # J32: 1966.3334 -> 2500
$ws.Cells.Item(32, 11).Value = 499  # K32: 760.4 -> 499
$ws.Cells.Item(32, 12).Value = 2500  # L32: 1966.3334 -> 2500
$ws.Cells.Item(32, 13).Value = -173  # M32: -434.4 -> -173
$ws.Cells.Item(32, 14).Value = -3152  # N32: -2618.3334 -> -3152
# Row 33
$ws.Cells.Item(33, 8).Value = 270.24243  # H33: 266.90323 -> 270.24243
$ws.Cells.Item(33, 9).Value = 203.44  # I33: 204.24 -> 203.44
$ws.Cells.Item(33, 10).Value = 479  # J33: 528 -> 479
$ws.Cells.Item(33, 11).Value = 203.44  # K33: 204.24 -> 203.44
$ws.Cells.Item(33, 12).Value = 479  # L33: 528 -> 479
$ws.Cells.Item(33, 13).Value = 25.56  # M33: 24.75999999999999 -> 25.56
$ws.Cells.Item(33, 14).Value = -937  # N33: -986 -> -937
# Row 34
$ws.Cells.Item(34, 8).Value = 5261  # H34: 3297 -> 5261
$ws.Cells.Item(34, 9).Value = 5261  # I34: 3297 -> 5261
$ws.Cells.Item(34, 11).Value = 5261  # K34: 3297 -> 5261
$ws.Cells.Item(34, 13).Value = -5058  # M34: -3094 -> -5058
# Row 36
$ws.Cells.Item(36, 8).Value = 5261  # H36: 3297 -> 5261
$ws.Cells.Item(36, 9).Value = 5261  # I36: 3297 -> 5261
$ws.Cells.Item(36, 11).Value = 5261  # K36: 3297 -> 5261
$ws.Cells.Item(36, 13).Value = -4546  # M36: -2582 -> -4546
# Row 40
$ws.Cells.Item(40, 8).Value = 2220.6  # H40: 2175.3333 -> 2220.6
$ws.Cells.Item(40, 9).Value = 2001  # I40: 1988 -> 2001
$ws.Cells.Item(40, 11).Value = 2001  # K40: 1988 -> 2001
$ws.Cells.Item(40, 13).Value = -1826  # M40: -1813 -> -1826
# Row 41
$ws.Cells.Item(41, 9).Value = 2736.375  # I41: 3098.7144 -> 2736.375
$ws.Cells.Item(41, 10).Value = 1643.75  # J41: 1483.3334 -> 1643.75
$ws.Cells.Item(41, 11).Value = 2736.375  # K41: 3098.7144 -> 2736.375
$ws.Cells.Item(41, 12).Value = 1643.75  # L41: 1483.3334 -> 1643.75
$ws.Cells.Item(41, 13).Value = -2296.375  # M41: -2658.7144 -> -2296.375
$ws.Cells.Item(41, 14).Value = -2523.75  # N41: -2363.3334 -> -2523.75
# Row 133
$ws.Cells.Item(133, 8).Value = 35251.125  # H133: 35594 -> 35251.125
$ws.Cells.Item(133, 10).Value = 35251.125  # J133: 35594 -> 35251.125
$ws.Cells.Item(133, 12).Value = 35251.125  # L133: 35594 -> 35251.125
$ws.Cells.Item(133, 14).Value = -45371.125  # N133: -45714 -> -45371.125
# Row 137
$ws.Cells.Item(137, 8).Value = 1379.4681  # H137: 1392.8 -> 1379.4681
$ws.Cells.Item(137, 9).Value = 1071  # I137: 1060.4762 -> 1071
$ws.Cells.Item(137, 10).Value = 1607.963  # J137: 1683.5834 -> 1607.963
$ws.Cells.Item(137, 11).Value = 3213  # K137: 3181.4286 -> 3213
$ws.Cells.Item(137, 12).Value = 4823.889  # L137: 5050.7502 -> 4823.889
$ws.Cells.Item(137, 13).Value = -663  # M137: -631.4286000000002 -> -663
$ws.Cells.Item(137, 14).Value = -9923.888999999999  # N137: -10150.7502 -> -9923.888999999999
# Row 138
$ws.Cells.Item(138, 8).Value = 1610.5238  # H138: 622633.9 -> 1610.5238
$ws.Cells.Item(138, 9).Value = 1100.5454  # I138: 998.41174 -> 1100.5454
$ws.Cells.Item(138, 10).Value = 1718.4038  # J138: 822026.4 -> 1718.4038
$ws.Cells.Item(138, 11).Value = 3301.6362  # K138: 2995.23522 -> 3301.6362
$ws.Cells.Item(138, 12).Value = 5155.2114  # L138: 2466079.2 -> 5155.2114
$ws.Cells.Item(138, 13).Value = 1838.3638  # M138: 2144.76478 -> 1838.3638
$ws.Cells.Item(138, 14).Value = -15435.2114  # N138: -2476359.2 -> -15435.2114

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 5245.4873  # H32: 5762.1143 -> 5245.4873
$ws.Cells.Item(32, 9).Value = 5330.8945  # I32: 5872.7646 -> 5330.8945
$ws.Cells.Item(32, 11).Value = 5330.8945  # K32: 5872.7646 -> 5330.8945
$ws.Cells.Item(32, 13).Value = -5043.8945  # M32: -5585.7646 -> -5043.8945
# Row 37
$ws.Cells.Item(37, 8).Value = 18900  # H37: 24360 -> 18900
$ws.Cells.Item(37, 10).Value = 0  # J37: 28000 -> 0
$ws.Cells.Item(37, 12).Value = 0  # L37: 28000 -> 0
$ws.Cells.Item(37, 14).ClearContents()  # N37: remove (was -28546)
# Row 44
$ws.Cells.Item(44, 8).Value = 13044  # H44: 23297.6 -> 13044
$ws.Cells.Item(44, 9).Value = 13044  # I44: 13544 -> 13044
$ws.Cells.Item(44, 10).Value = 0  # J44: 29800 -> 0
$ws.Cells.Item(44, 11).Value = 13044  # K44: 13544 -> 13044
$ws.Cells.Item(44, 12).Value = 0  # L44: 29800 -> 0
$ws.Cells.Item(44, 13).Value = -12556  # M44: -13056 -> -12556
$ws.Cells.Item(44, 14).ClearContents()  # N44: remove (was -30776)
# Row 45
$ws.Cells.Item(45, 8).Value = 1298.7693  # H45: 1241.7142 -> 1298.7693
$ws.Cells.Item(45, 9).Value = 1207.7778  # I45: 1137 -> 1207.7778
$ws.Cells.Item(45, 11).Value = 1207.7778  # K45: 1137 -> 1207.7778
$ws.Cells.Item(45, 13).Value = -830.7778000000001  # M45: -760 -> -830.7778000000001
# Row 55
$ws.Cells.Item(55, 8).Value = 37800  # H55: 37966.332 -> 37800
$ws.Cells.Item(55, 10).Value = 37800  # J55: 37966.332 -> 37800
$ws.Cells.Item(55, 12).Value = 37800  # L55: 37966.332 -> 37800
$ws.Cells.Item(55, 14).Value = -38430  # N55: -38596.332 -> -38430
# Row 122
$ws.Cells.Item(122, 8).Value = 2312.4  # H122: 2303 -> 2312.4
$ws.Cells.Item(122, 9).Value = 1638.75  # I122: 1915.7142 -> 1638.75
$ws.Cells.Item(122, 10).Value = 5007  # J122: 5014 -> 5007
$ws.Cells.Item(122, 11).Value = 4916.25  # K122: 5747.142599999999 -> 4916.25
$ws.Cells.Item(122, 12).Value = 15021  # L122: 15042 -> 15021
$ws.Cells.Item(122, 13).Value = -2466.25  # M122: -3297.142599999999 -> -2466.25
$ws.Cells.Item(122, 14).Value = -19921  # N122: -19942 -> -19921
# Row 139
$ws.Cells.Item(139, 8).Value = 27348.572  # H139: 29370 -> 27348.572
$ws.Cells.Item(139, 10).Value = 27348.572  # J139: 29370 -> 27348.572
$ws.Cells.Item(139, 12).Value = 27348.572  # L139: 29370 -> 27348.572
$ws.Cells.Item(139, 14).Value = -37628.572  # N139: -39650 -> -37628.572

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Cells.Item(105, 8).Value = 336630300  # H105: 168316180 -> 336630300
$ws.Cells.Item(105, 9).Value = 336630300  # I105: 201978980 -> 336630300
$ws.Cells.Item(105, 10).Value = 0  # J105: 2200 -> 0
$ws.Cells.Item(105, 11).Value = 336630300  # K105: 201978980 -> 336630300
$ws.Cells.Item(105, 12).Value = 0  # L105: 2200 -> 0
$ws.Cells.Item(105, 13).Value = -336628553  # M105: -201977233 -> -336628553
$ws.Cells.Item(105, 14).ClearContents()  # N105: remove (was -5694)

$ws = $wb.Worksheets.Item("CRP")
# Row 134
$ws.Cells.Item(134, 8).Value = 20835282  # H134: 21741186 -> 20835282
$ws.Cells.Item(134, 9).Value = 2038.1  # I134: 2256.647 -> 2038.1
$ws.Cells.Item(134, 10).Value = 125001500  # J134: 83334820 -> 125001500
$ws.Cells.Item(134, 11).Value = 6114.299999999999  # K134: 6769.941 -> 6114.299999999999
$ws.Cells.Item(134, 12).Value = 375004500  # L134: 250004460 -> 375004500
$ws.Cells.Item(134, 13).Value = -3579.299999999999  # M134: -4234.941 -> -3579.299999999999
$ws.Cells.Item(134, 14).Value = -375009570  # N134: -250009530 -> -375009570

$ws = $wb.Worksheets.Item("CUL")
# Row 55
$ws.Cells.Item(55, 8).Value = 3126.4285  # H55: 3183.3333 -> 3126.4285
$ws.Cells.Item(55, 10).Value = 3126.4285  # J55: 3183.3333 -> 3126.4285
$ws.Cells.Item(55, 12).Value = 9379.2855  # L55: 9549.999899999999 -> 9379.2855
$ws.Cells.Item(55, 14).Value = -9733.2855  # N55: -9903.999899999999 -> -9733.2855
# Row 132
$ws.Cells.Item(132, 8).Value = 804.3  # H132: 698.7692 -> 804.3
$ws.Cells.Item(132, 9).Value = 830.5  # I132: 698.7692 -> 830.5
$ws.Cells.Item(132, 10).Value = 699.5  # J132: 0 -> 699.5
$ws.Cells.Item(132, 11).Value = 7474.5  # K132: 6288.922799999999 -> 7474.5
$ws.Cells.Item(132, 12).Value = 6295.5  # L132: 0 -> 6295.5
$ws.Cells.Item(132, 13).Value = -4944.5  # M132: -3758.922799999999 -> -4944.5
$ws.Cells.Item(132, 14).Value = -11355.5  # N132: <MISSING> -> -11355.5

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Cells.Item(122, 8).Value = 6580312.5  # H122: 6758151 -> 6580312.5
$ws.Cells.Item(122, 9).Value = 1342.3334  # I122: 1382.6154 -> 1342.3334
$ws.Cells.Item(122, 11).Value = 4027.0002  # K122: 4147.8462 -> 4027.0002
$ws.Cells.Item(122, 13).Value = -1577.0002  # M122: -1697.8462 -> -1577.0002
# Row 126
$ws.Cells.Item(126, 8).Value = 1994.6  # H126: 2167.3333 -> 1994.6
$ws.Cells.Item(126, 9).Value = 1670.5  # I126: 1809.1111 -> 1670.5
$ws.Cells.Item(126, 10).Value = 2642.8  # J126: 2704.6667 -> 2642.8
$ws.Cells.Item(126, 11).Value = 5011.5  # K126: 5427.3333 -> 5011.5
$ws.Cells.Item(126, 12).Value = 7928.400000000001  # L126: 8114.000100000001 -> 7928.400000000001
$ws.Cells.Item(126, 13).Value = -2541.5  # M126: -2957.3333 -> -2541.5
$ws.Cells.Item(126, 14).Value = -12868.4  # N126: -13054.0001 -> -12868.4

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 702.75  # H22: 700.1539 -> 702.75
$ws.Cells.Item(22, 9).Value = 464.5  # I22: 476.16666 -> 464.5
$ws.Cells.Item(22, 10).Value = 941  # J22: 892.1429000000001 -> 941
$ws.Cells.Item(22, 11).Value = 464.5  # K22: 476.16666 -> 464.5
$ws.Cells.Item(22, 12).Value = 941  # L22: 892.1429000000001 -> 941
$ws.Cells.Item(22, 13).Value = -169.5  # M22: -181.16666 -> -169.5
$ws.Cells.Item(22, 14).Value = -1531  # N22: -1482.1429 -> -1531
# Row 27
$ws.Cells.Item(27, 8).Value = 702.75  # H27: 700.1539 -> 702.75
$ws.Cells.Item(27, 9).Value = 464.5  # I27: 476.16666 -> 464.5
$ws.Cells.Item(27, 10).Value = 941  # J27: 892.1429000000001 -> 941
$ws.Cells.Item(27, 11).Value = 464.5  # K27: 476.16666 -> 464.5
$ws.Cells.Item(27, 12).Value = 941  # L27: 892.1429000000001 -> 941
$ws.Cells.Item(27, 13).Value = -357.5  # M27: -369.16666 -> -357.5
$ws.Cells.Item(27, 14).Value = -1155  # N27: -1106.1429 -> -1155
# Row 46
$ws.Cells.Item(46, 8).Value = 3588.6428  # H46: 3648.7144 -> 3588.6428
$ws.Cells.Item(46, 9).Value = 534.1429000000001  # I46: 596.6667 -> 534.1429000000001
$ws.Cells.Item(46, 10).Value = 6643.143  # J46: 5937.75 -> 6643.143
$ws.Cells.Item(46, 11).Value = 534.1429000000001  # K46: 596.6667 -> 534.1429000000001
$ws.Cells.Item(46, 12).Value = 6643.143  # L46: 5937.75 -> 6643.143
$ws.Cells.Item(46, 13).Value = -346.1429000000001  # M46: -408.6667 -> -346.1429000000001
$ws.Cells.Item(46, 14).Value = -7019.143  # N46: -6313.75 -> -7019.143
# Row 55
$ws.Cells.Item(55, 8).Value = 313.55554  # H55: 299.6842 -> 313.55554
$ws.Cells.Item(55, 10).Value = 383.14285  # J55: 341.5 -> 383.14285
$ws.Cells.Item(55, 12).Value = 383.14285  # L55: 341.5 -> 383.14285
$ws.Cells.Item(55, 14).Value = -729.14285  # N55: -687.5 -> -729.14285
# Row 100
$ws.Cells.Item(100, 8).Value = 1386.75  # H100: 1389 -> 1386.75
$ws.Cells.Item(100, 9).Value = 1386.75  # I100: 1389 -> 1386.75
$ws.Cells.Item(100, 11).Value = 1386.75  # K100: 1389 -> 1386.75
$ws.Cells.Item(100, 13).Value = -845.75  # M100: -848 -> -845.75
# Row 122
$ws.Cells.Item(122, 8).Value = 250000000  # H122: 62501696 -> 250000000
$ws.Cells.Item(122, 9).Value = 250000000  # I122: 83334800 -> 250000000
$ws.Cells.Item(122, 10).Value = 0  # J122: 2400 -> 0
$ws.Cells.Item(122, 11).Value = 750000000  # K122: 250004400 -> 750000000
$ws.Cells.Item(122, 12).Value = 0  # L122: 7200 -> 0
$ws.Cells.Item(122, 13).Value = -749997550  # M122: -250001950 -> -749997550
$ws.Cells.Item(122, 14).ClearContents()  # N122: remove (was -12100)
# Row 134
$ws.Cells.Item(134, 8).Value = 34438.4  # H134: 34981.668 -> 34438.4
$ws.Cells.Item(134, 10).Value = 34438.4  # J134: 34981.668 -> 34438.4
$ws.Cells.Item(134, 12).Value = 34438.4  # L134: 34981.668 -> 34438.4
$ws.Cells.Item(134, 14).Value = -44578.4  # N134: -45121.668 -> -44578.4

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Cells.Item(122, 8).Value = 16667906  # H122: 17858442 -> 16667906
$ws.Cells.Item(122, 9).Value = 16667906  # I122: 17858442 -> 16667906
$ws.Cells.Item(122, 11).Value = 50003718  # K122: 53575326 -> 50003718
$ws.Cells.Item(122, 13).Value = -50001268  # M122: -53572876 -> -50001268
# Row 124
$ws.Cells.Item(124, 8).Value = 35333.332  # H124: 20173.375 -> 35333.332
$ws.Cells.Item(124, 10).Value = 35333.332  # J124: 20173.375 -> 35333.332
$ws.Cells.Item(124, 12).Value = 35333.332  # L124: 20173.375 -> 35333.332
$ws.Cells.Item(124, 14).Value = -45153.332  # N124: -29993.375 -> -45153.332
# Row 133
$ws.Cells.Item(133, 8).Value = 33900  # H133: 37538.332 -> 33900
$ws.Cells.Item(133, 10).Value = 33900  # J133: 37538.332 -> 33900
$ws.Cells.Item(133, 12).Value = 33900  # L133: 37538.332 -> 33900
$ws.Cells.Item(133, 14).Value = -44020  # N133: -47658.332 -> -44020
# Row 135
$ws.Cells.Item(135, 8).Value = 715000  # H135: 88343.2 -> 715000
$ws.Cells.Item(135, 10).Value = 715000  # J135: 88343.2 -> 715000
$ws.Cells.Item(135, 12).Value = 715000  # L135: 88343.2 -> 715000
$ws.Cells.Item(135, 14).Value = -725140  # N135: -98483.2 -> -725140
